$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings formatted like "1.853.27" (thousand-separator style).
# Force the whole column to Text so numeric-looking updates ("1.003", "233.02", ...)
# are stored as text instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.396.65"
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").Value = "1.849.92"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "233.02"
$ws.Range("E5").Value = "  +1.75%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "0.4737"
$ws.Range("E7").Value = "  +2.83%  "
$ws.Range("D8").Value = "0.2741"
$ws.Range("E8").Value = "  +2.17%  "
$ws.Range("D9").Value = "0.06305"
$ws.Range("E9").Value = "  +1.93%  "
$ws.Range("D10").Value = "17.57"
$ws.Range("E10").Value = "  +10.59%  "
$ws.Range("D11").Value = "1.863.11"
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("D12").Value = "0.07445"
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("D13").Value = "4.950"
$ws.Range("E13").Value = "  +1.69%  "
$ws.Range("D14").Value = "84.44"
$ws.Range("E14").Value = "  +2.51%  "
$ws.Range("D15").Value = "0.6233"
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("D16").Value = "30.401.18"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("D17").Value = "246.29"
$ws.Range("E17").Value = "  +9.84%  "
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "12.67"
$ws.Range("E19").Value = "  +3.63%  "
$ws.Range("D20").Value = "0.000007302"
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").Value = "4.913"
$ws.Range("E22").Value = "  +2.38%  "
$ws.Range("D23").Value = "5.887"
$ws.Range("E23").Value = "  +1.21%  "
$ws.Range("D24").Value = "9.104"
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("D25").Value = "164.56"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "17.93"
$ws.Range("E26").Value = "  +2.32%  "
$ws.Range("D27").Value = "1.866"
$ws.Range("E27").Value = "  +1.89%  "
$ws.Range("D29").Value = "1.355"
$ws.Range("E29").Value = "  -1.28%  "
$ws.Range("D30").Value = "4.029"
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("D31").Value = "3.819"
$ws.Range("E31").Value = "  +2.21%  "
$ws.Range("D32").Value = "0.04836"
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("D33").Value = "1.126"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").Value = "0.6947"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").Value = "2.702"
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("D36").Value = "0.01889"
$ws.Range("E36").Value = "  +5.67%  "
$ws.Range("D38").Value = "1.995"
$ws.Range("E38").Value = "  +5.12%  "
$ws.Range("D39").Value = "0.8736"
$ws.Range("E39").Value = "  -0.79%  "
$ws.Range("D40").Value = "106.16"
$ws.Range("E40").Value = "  +3.37%  "
$ws.Range("D42").Value = "5.521"
$ws.Range("E42").Value = "  +1.86%  "
$ws.Range("D43").Value = "0.4042"
$ws.Range("E43").Value = "  +1.92%  "
$ws.Range("D44").Value = "7.161"
$ws.Range("E44").Value = "  +4.52%  "
$ws.Range("D45").Value = "62.85"
$ws.Range("E45").Value = "  +6.81%  "
$ws.Range("D46").Value = "0.1194"
$ws.Range("E46").Value = "  +1.47%  "
$ws.Range("D47").Value = "33.61"
$ws.Range("E47").Value = "  +3.98%  "
$ws.Range("D48").Value = "8.544"
$ws.Range("E48").Value = "  +2.23%  "
$ws.Range("D49").Value = "0.05521"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").Value = "1.345"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").Value = "0.3671"
$ws.Range("E51").Value = "  +2.31%  "

# Rows where only the Volume(1h) column changed
$ws.Range("E28").Value = "  +1.48%  "
$ws.Range("E41").Value = "  +0.64%  "
